# Commit: Wed, May 06, 2020 10:06:40 AM
#
# Two logical changes in the source OOXML:
#   1. The three tables on slides 14-16 switch their table style from
#      {8B5209A2-B309-4E1D-AE42-B68ACD26F32C} to
#      {69A180DD-6C25-46BF-8617-8766DBD19B8B}.
#   2. ppt/theme/theme1.xml and ppt/theme/theme2.xml swap their full
#      contents (the "Integral"/"Red Violet" theme and the
#      "Office Theme"/"Office" theme trade places). Since theme1.xml is
#      the theme actually driving the slide master / slides, and the
#      only structural differences between the two themes are the
#      <a:clrScheme> color values (font/format schemes are identical),
#      that swap is reproduced here by rewriting the 12 theme colors
#      that theme1.xml exposes through the presentation's
#      ThemeColorScheme object to the target ("Office") palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style every table in the deck that currently uses the old
#        style id. ---------------------------------------------------
$oldStyleId = "{8B5209A2-B309-4E1D-AE42-B68ACD26F32C}"
$newStyleId = "{69A180DD-6C25-46BF-8617-8766DBD19B8B}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Swap the theme colours (theme1.xml currently holds the
#        "Integral" / "Red Violet" palette; the target state is the
#        stock "Office" palette that used to live in theme2.xml). -----
$tcs = $p.Slides.Item(1).ThemeColorScheme

function Set-ThemeColor([int]$index, [int]$r, [int]$g, [int]$b) {
    $tcs.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-ThemeColor 1  0x00 0x00 0x00   # dk1
Set-ThemeColor 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor 12 0x95 0x4F 0x72   # folHlink
